$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 3847.5
$ws.Range("J29").Value = 4996.6665
$ws.Range("L29").Value = 14989.9995
$ws.Range("N29").Value = -15551.9995

$ws.Range("H69").Value = 15045.25
$ws.Range("I69").Value = 15090.5
$ws.Range("K69").Value = 45271.5
$ws.Range("M69").Value = -44397.5

$ws.Range("H72").Value = 15045.25
$ws.Range("I72").Value = 15090.5
$ws.Range("K72").Value = 135814.5
$ws.Range("M72").Value = -131446.5

$ws.Range("H106").Value = 16500
$ws.Range("I106").Value = 2400
$ws.Range("K106").Value = 2400
$ws.Range("M106").Value = -1769

$ws.Range("H112").Value = 2651.6667
$ws.Range("J112").Value = 2473.75
$ws.Range("L112").Value = 7421.25
$ws.Range("N112").Value = -9637.25

$ws.Range("H138").Value = 11262.058
$ws.Range("I138").Value = 12835.272
$ws.Range("J138").Value = 10541
$ws.Range("K138").Value = 38505.81600000001
$ws.Range("L138").Value = 31623
$ws.Range("M138").Value = -33365.81600000001
$ws.Range("N138").Value = -41903

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 720.99
$ws.Range("I32").Value = 720.99
$ws.Range("K32").Value = 720.99
$ws.Range("M32").Value = -433.99

$ws.Range("H61").Value = 7694310.5
$ws.Range("I61").Value = 1916.4
$ws.Range("J61").Value = 33335624
$ws.Range("K61").Value = 1916.4
$ws.Range("L61").Value = 33335624
$ws.Range("M61").Value = -1704.4
$ws.Range("N61").Value = -33336048

$ws.Range("H74").Value = 1596668.6
$ws.Range("I74").Value = 3091961.5
$ws.Range("K74").Value = 3091961.5
$ws.Range("M74").Value = -3091087.5

$ws.Range("H77").Value = 1596668.6
$ws.Range("I77").Value = 3091961.5
$ws.Range("K77").Value = 15459807.5
$ws.Range("M77").Value = -15455439.5

$ws.Range("H136").Value = 7694310.5
$ws.Range("I136").Value = 1916.4
$ws.Range("J136").Value = 33335624
$ws.Range("K136").Value = 5749.200000000001
$ws.Range("L136").Value = 100006872
$ws.Range("M136").Value = -3199.200000000001
$ws.Range("N136").Value = -100011972

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 32001.5
$ws.Range("J103").Value = 32001.5
$ws.Range("L103").Value = 32001.5
$ws.Range("N103").Value = -34345.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 92703.91
$ws.Range("I16").Value = 1804.7778
$ws.Range("K16").Value = 1804.7778
$ws.Range("M16").Value = -1517.7778

$ws.Range("H86").Value = 17351.625
$ws.Range("I86").Value = 16437.666
$ws.Range("J86").Value = 18874.889
$ws.Range("K86").Value = 16437.666
$ws.Range("L86").Value = 18874.889
$ws.Range("M86").Value = -15314.666
$ws.Range("N86").Value = -21120.889

$ws.Range("H89").Value = 17351.625
$ws.Range("I89").Value = 16437.666
$ws.Range("J89").Value = 18874.889
$ws.Range("K89").Value = 82188.33
$ws.Range("L89").Value = 94374.44499999999
$ws.Range("M89").Value = -76572.33
$ws.Range("N89").Value = -105606.445

$ws.Range("H113").Value = 92703.91
$ws.Range("I113").Value = 1804.7778
$ws.Range("K113").Value = 1804.7778
$ws.Range("M113").Value = 365.2221999999999

$ws.Range("H132").Value = 5200.0415
$ws.Range("I132").Value = 4926.143
$ws.Range("K132").Value = 14778.429
$ws.Range("M132").Value = -12248.429

$ws.Range("H134").Value = 1386.0454
$ws.Range("I134").Value = 1314
$ws.Range("J134").Value = 2899
$ws.Range("K134").Value = 3942
$ws.Range("L134").Value = 8697
$ws.Range("M134").Value = -1407
$ws.Range("N134").Value = -13767

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2776369.2
$ws.Range("I4").Value = 3101529
$ws.Range("K4").Value = 9304587
$ws.Range("M4").Value = -9304475

$ws.Range("H59").Value = 1875
$ws.Range("I59").Value = 1833.3334
$ws.Range("K59").Value = 5500.0002
$ws.Range("M59").Value = -4960.0002

$ws.Range("H87").Value = 11998.5
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H90").Value = 11998.5
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H129").Value = 460136.1
$ws.Range("I129").Value = 1430929.4
$ws.Range("J129").Value = 7099.2
$ws.Range("K129").Value = 4292788.199999999
$ws.Range("L129").Value = 21297.6
$ws.Range("M129").Value = -4287788.199999999
$ws.Range("N129").Value = -31297.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4926.8887
$ws.Range("I122").Value = 5115.4165
$ws.Range("J122").Value = 4549.8335
$ws.Range("K122").Value = 15346.2495
$ws.Range("L122").Value = 13649.5005
$ws.Range("M122").Value = -12896.2495
$ws.Range("N122").Value = -18549.5005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5032.5
$ws.Range("I16").Value = 1076.625
$ws.Range("J16").Value = 12944.25
$ws.Range("K16").Value = 1076.625
$ws.Range("L16").Value = 12944.25
$ws.Range("M16").Value = -906.625
$ws.Range("N16").Value = -13284.25

$ws.Range("H68").Value = 12062.9
$ws.Range("I68").Value = 9909.25
$ws.Range("J68").Value = 16370.2
$ws.Range("K68").Value = 9909.25
$ws.Range("L68").Value = 16370.2
$ws.Range("M68").Value = -9160.25
$ws.Range("N68").Value = -17868.2

$ws.Range("H71").Value = 12062.9
$ws.Range("I71").Value = 9909.25
$ws.Range("J71").Value = 16370.2
$ws.Range("K71").Value = 49546.25
$ws.Range("L71").Value = 81851
$ws.Range("M71").Value = -45802.25
$ws.Range("N71").Value = -89339

$ws.Range("H93").Value = 3994.2222
$ws.Range("I93").Value = 2389.8
$ws.Range("K93").Value = 2389.8
$ws.Range("M93").Value = -1141.8

$ws.Range("H136").Value = 6759195.5
$ws.Range("I136").Value = 4034689.8
$ws.Range("K136").Value = 12104069.4
$ws.Range("M136").Value = -12101519.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 2999.5
$ws.Range("J6").Value = 2999.5
$ws.Range("L6").Value = 2999.5
$ws.Range("N6").Value = -3229.5

$ws.Range("H81").Value = 76926104
$ws.Range("I81").Value = 2835.1
$ws.Range("J81").Value = 333337000
$ws.Range("K81").Value = 5670.2
$ws.Range("L81").Value = 666674000
$ws.Range("M81").Value = -4609.2
$ws.Range("N81").Value = -666676122

$ws.Range("H84").Value = 76926104
$ws.Range("I84").Value = 2835.1
$ws.Range("J84").Value = 333337000
$ws.Range("K84").Value = 28351
$ws.Range("L84").Value = 3333370000
$ws.Range("M84").Value = -23047
$ws.Range("N84").Value = -3333380608

$ws.Range("H107").Value = 40003316
$ws.Range("I107").Value = 76924936
$ws.Range("J107").Value = 4891.8335
$ws.Range("K107").Value = 230774808
$ws.Range("L107").Value = 14675.5005
$ws.Range("M107").Value = -230772888
$ws.Range("N107").Value = -18515.5005
